# Generate Report for Handback
# Updates the "Latest Target File" (I), "Latest Target Datetime" (J... actually
# Latest Target File display uses column J for the handback xlf filename),
# "Latest Handback DateTime" (K) and "Error Detail" (P) columns for the
# 75413ddf-d041-4eae-903a-3dc4ab6e1231 row on both the zh-cn and de-de sheets,
# reflecting that a newer handback was produced but is out of date versus the
# latest handoff.

$wb = $excel.ActiveWorkbook

$targetMdFile = "75413ddf-d041-4eae-903a-3dc4ab6e1231.md"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/63e4ea37be9a4a83ccdad06d0557b131db98708e/e2e/75413ddf-d041-4eae-903a-3dc4ab6e1231.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7c47c641d10e7f4962ce06b60563a39fb3cbc339/e2e/75413ddf-d041-4eae-903a-3dc4ab6e1231.md."

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("J7").Value = "75413ddf-d041-4eae-903a-3dc4ab6e1231.e5b04453cd0f64db7bb46fe7a97e7314a259e8e2.zh-cn.xlf"
$wsZh.Range("K7").Value = "2016-09-03 17:02:44"
$wsZh.Range("P7").Value = $errorDetail

$wsZh.Hyperlinks.Add($wsZh.Range("I7"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/63e4ea37be9a4a83ccdad06d0557b131db98708e/e2e/75413ddf-d041-4eae-903a-3dc4ab6e1231.md", [Type]::Missing, [Type]::Missing, $targetMdFile)

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("J7").Value = "75413ddf-d041-4eae-903a-3dc4ab6e1231.e5b04453cd0f64db7bb46fe7a97e7314a259e8e2.de-de.xlf"
$wsDe.Range("K7").Value = "2016-09-03 17:02:51"
$wsDe.Range("P7").Value = $errorDetail

$wsDe.Hyperlinks.Add($wsDe.Range("I7"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/63e4ea37be9a4a83ccdad06d0557b131db98708e/e2e/75413ddf-d041-4eae-903a-3dc4ab6e1231.md", [Type]::Missing, [Type]::Missing, $targetMdFile)
